$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.379369
$ws.Range("H2").Value = 31.138107
$ws.Range("I2").Value = 0.01614698522449884
$ws.Range("J2").Value = 0.01614698522449883
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.764494666666667
$ws.Range("N2").Value = 5.293483999999999
$ws.Range("O2").Value = 0.0152564507897189
$ws.Range("P2").Value = 0.0152564507897189
$ws.Range("Q2").Value = 18.31434124386533
$ws.Range("R2").Value = 164.829071194788
$ws.Range("S2").Value = 0.0002463456854798847
$ws.Range("T2").Value = 0.0002463456854798846

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.379369
$ws.Range("H3").Value = 31.138107
$ws.Range("I3").Value = 0.01614698522449884
$ws.Range("J3").Value = 0.01614698522449883
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 93.97803500000002
$ws.Range("N3").Value = 281.934105
$ws.Range("O3").Value = 0.8125676395500474
$ws.Range("P3").Value = 0.8125676395500474
$ws.Range("Q3").Value = 975.4327031599153
$ws.Range("R3").Value = 8778.894328439237
$ws.Range("S3").Value = 0.01312051766972051
$ws.Range("T3").Value = 0.01312051766972051

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.379369
$ws.Range("H4").Value = 31.138107
$ws.Range("I4").Value = 0.01614698522449884
$ws.Range("J4").Value = 0.01614698522449883
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 19.913116
$ws.Range("N4").Value = 59.73934800000001
$ws.Range("O4").Value = 0.1721759096602337
$ws.Range("P4").Value = 0.1721759096602337
$ws.Range("Q4").Value = 206.685578903804
$ws.Range("R4").Value = 1860.170210134236
$ws.Range("S4").Value = 0.00278012186929844
$ws.Range("T4").Value = 0.00278012186929844

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 604.0312093333333
$ws.Range("H5").Value = 1812.093628
$ws.Range("I5").Value = 0.9396797639857967
$ws.Range("J5").Value = 0.9396797639857967
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.764494666666667
$ws.Range("N5").Value = 5.293483999999999
$ws.Range("O5").Value = 0.0152564507897189
$ws.Range("P5").Value = 0.0152564507897189
$ws.Range("Q5").Value = 1065.809847368884
$ws.Range("R5").Value = 9592.288626319951
$ws.Range("S5").Value = 0.01433617807734398
$ws.Range("T5").Value = 0.01433617807734398

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 604.0312093333333
$ws.Range("H6").Value = 1812.093628
$ws.Range("I6").Value = 0.9396797639857967
$ws.Range("J6").Value = 0.9396797639857967
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 93.97803500000002
$ws.Range("N6").Value = 281.934105
$ws.Range("O6").Value = 0.8125676395500474
$ws.Range("P6").Value = 0.8125676395500474
$ws.Range("Q6").Value = 56765.66613182034
$ws.Range("R6").Value = 510890.995186383
$ws.Range("S6").Value = 0.7635533677548845
$ws.Range("T6").Value = 0.7635533677548845

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 604.0312093333333
$ws.Range("H7").Value = 1812.093628
$ws.Range("I7").Value = 0.9396797639857967
$ws.Range("J7").Value = 0.9396797639857967
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 19.913116
$ws.Range("N7").Value = 59.73934800000001
$ws.Range("O7").Value = 0.1721759096602337
$ws.Range("P7").Value = 0.1721759096602337
$ws.Range("Q7").Value = 12028.14353907495
$ws.Range("R7").Value = 108253.2918516746
$ws.Range("S7").Value = 0.1617902181535683
$ws.Range("T7").Value = 0.1617902181535683

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 28.39480333333333
$ws.Range("H8").Value = 85.18441
$ws.Range("I8").Value = 0.04417325078970442
$ws.Range("J8").Value = 0.04417325078970442
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.764494666666667
$ws.Range("N8").Value = 5.293483999999999
$ws.Range("O8").Value = 0.0152564507897189
$ws.Range("P8").Value = 0.0152564507897189
$ws.Range("Q8").Value = 50.10247904271555
$ws.Range("R8").Value = 450.9223113844399
$ws.Range("S8").Value = 0.0006739270268950371
$ws.Range("T8").Value = 0.0006739270268950371

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 28.39480333333333
$ws.Range("H9").Value = 85.18441
$ws.Range("I9").Value = 0.04417325078970442
$ws.Range("J9").Value = 0.04417325078970442
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 93.97803500000002
$ws.Range("N9").Value = 281.934105
$ws.Range("O9").Value = 0.8125676395500474
$ws.Range("P9").Value = 0.8125676395500474
$ws.Range("Q9").Value = 2668.487821478117
$ws.Range("R9").Value = 24016.39039330305
$ws.Range("S9").Value = 0.03589375412544239
$ws.Range("T9").Value = 0.03589375412544239

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 28.39480333333333
$ws.Range("H10").Value = 85.18441
$ws.Range("I10").Value = 0.04417325078970442
$ws.Range("J10").Value = 0.04417325078970442
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 19.913116
$ws.Range("N10").Value = 59.73934800000001
$ws.Range("O10").Value = 0.1721759096602337
$ws.Range("P10").Value = 0.1721759096602337
$ws.Range("Q10").Value = 565.4290125738534
$ws.Range("R10").Value = 5088.86111316468
$ws.Range("S10").Value = 0.007605569637366994
$ws.Range("T10").Value = 0.007605569637366996

